$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Name -> "جيانا", file_date -> "01-01", path formula includes both Name and file_date
$ws.Range("A2").Value = "جيانا"
$ws.Range("C2").Value = "01-01"
$ws.Range("B2").Formula = "=""C:\Users\AL-Thuraya\Documents\Projects\OPOST_AUTOMATION\samples\""&A2&"" ""&C2"

# Row 3: Name -> "رند", file_date -> "05-25", path formula includes both Name and file_date
$ws.Range("A3").Value = "رند"
$ws.Range("C3").Value = "05-25"
$ws.Range("B3").Formula = "=""C:\Users\AL-Thuraya\Documents\Projects\OPOST_AUTOMATION\samples\""&A3&"" ""&C3"

# Update the saved selection to C8
$ws.Range("C8").Select()
